$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lernjournal Aufgabe 2")
Write-Host $ws.Name
Write-Host $ws.Range("A31").Value
